# Add the "Wireframes" peer-review row (review comments on the wireframes,
# done by Nesma Bahgat / reviewed by Asmaa Hamdy) to the TAWA peer review
# sheet, right below the existing "HLD (Design Document)" row (row 11).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New row 12: reuse the existing formatting already present in the sheet
#     so no redundant style records are introduced.
#       - C11 already carries the bordered / wrapped / vertically-centered
#         look that A12:C12 need.
#       - F11 already carries the date look that F12 needs.
#       - E4 already carries the wrapped "comments" look that E12 needs.
$ws.Range("C11").Copy() | Out-Null
$ws.Range("A12:C12").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

$ws.Range("F11").Copy() | Out-Null
$ws.Range("F12").PasteSpecial(-4122) | Out-Null        # xlPasteFormats

$ws.Range("E4").Copy() | Out-Null
$ws.Range("E12").PasteSpecial(-4122) | Out-Null        # xlPasteFormats

$excel.CutCopyMode = 0

# --- Values for the new "Wireframes" review row ---
$ws.Range("A12").Value = "Wireframes"
$ws.Range("B12").Value = "Nesma Bahgat"
$ws.Range("C12").Value = "Asmaa Hamdy"
$ws.Range("E12").Value = "1- Add show details button to each trip`n2- The content of the side-bars doesn't exist`nThe following screens don't exist:-`n3- What will the user see after submitting the flight?`n4- Where will the user be directed to after signing-up?"
$ws.Range("F12").Value = 43743

$ws.Rows.Item(12).RowHeight = 110.25

# --- Leave the sheet with the same selection the author ended up with ---
$ws.Range("E15").Select() | Out-Null
